$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 3 new rows (140, 141, 142) continuing the Day_Number / Date series
$ws.Range("A140").Value = 139
$ws.Range("B140").Value = 43604

$ws.Range("A141").Value = 140
$ws.Range("B141").Value = 43605

$ws.Range("A142").Value = 141
$ws.Range("B142").Value = 43606

# Match formatting used by the preceding data row (A: style 3, B: style 4)
$ws.Range("A139:B139").Copy()
$ws.Range("A140:B142").PasteSpecial(-4122)

# Update selection to mirror the extended range (A141:A142 active)
$ws.Range("A141:A142").Activate()
